$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 31021
$ws.Range("E2").Value = 1670
$ws.Range("F2").Value = 1670
$ws.Range("G2").Value = 1627
$ws.Range("H2").Value = 2366
$ws.Range("I2").Value = 2295
$ws.Range("J2").Value = 71
$ws.Range("K2").Value = 257747
$ws.Range("L2").Value = 223479
$ws.Range("M2").Value = 34268
$ws.Range("N2").Value = 34268
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 3942
$ws.Range("Q2").Value = -15172
$ws.Range("R2").Value = -2950
$ws.Range("S2").Value = 20271
$ws.Range("T2").Value = 216
$ws.Range("V2").Value = 5927
$ws.Range("W2").Value = 5.38
$ws.Range("X2").Value = 7.63
$ws.Range("Y2").Value = 6.83
$ws.Range("Z2").Value = 0.99
$ws.Range("AA2").Value = 652.14
$ws.Range("AB2").Value = 815.73
$ws.Range("AC2").Value = 2911
$ws.Range("AD2").Value = 14.84
$ws.Range("AE2").Value = 45541
$ws.Range("AF2").Value = 0.95
$ws.Range("AG2").Value = 630
$ws.Range("AH2").Value = 1.46
$ws.Range("AI2").Value = 20.65
$ws.Range("AJ2").Value = 78822179
$ws.Range("U2").ClearContents()

# Row 3
$ws.Range("D3").Value = 39405
$ws.Range("E3").Value = 3767
$ws.Range("F3").Value = 3767
$ws.Range("G3").Value = 3649
$ws.Range("H3").Value = 2750
$ws.Range("I3").Value = 2750
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 309944
$ws.Range("L3").Value = 274706
$ws.Range("M3").Value = 35238
$ws.Range("N3").Value = 35238
$ws.Range("P3").Value = 3942
$ws.Range("Q3").Value = -14172
$ws.Range("R3").Value = -8397
$ws.Range("S3").Value = 23933
$ws.Range("T3").Value = 170
$ws.Range("V3").Value = 10639
$ws.Range("W3").Value = 9.56
$ws.Range("X3").Value = 6.98
$ws.Range("Y3").Value = 7.91
$ws.Range("Z3").Value = 0.97
$ws.Range("AA3").Value = 779.5700000000001
$ws.Range("AB3").Value = 872
$ws.Range("AC3").Value = 3489
$ws.Range("AD3").Value = 11.61
$ws.Range("AE3").Value = 48609
$ws.Range("AF3").Value = 0.83
$ws.Range("AG3").Value = 970
$ws.Range("AH3").Value = 2.39
$ws.Range("AI3").Value = 25.49
$ws.Range("AJ3").Value = 78822179
$ws.Range("O3").ClearContents()
$ws.Range("U3").ClearContents()

# Row 4
$ws.Range("D4").Value = 44285
$ws.Range("E4").Value = 2117
$ws.Range("F4").Value = 2117
$ws.Range("G4").Value = 2299
$ws.Range("H4").Value = 1742
$ws.Range("I4").Value = 1742
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 332299
$ws.Range("L4").Value = 294030
$ws.Range("M4").Value = 38270
$ws.Range("N4").Value = 38269
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 3942
$ws.Range("Q4").Value = 14669
$ws.Range("R4").Value = -25339
$ws.Range("S4").Value = 15222
$ws.Range("T4").Value = 192
$ws.Range("V4").Value = 10083
$ws.Range("W4").Value = 4.78
$ws.Range("X4").Value = 3.94
$ws.Range("Y4").Value = 4.74
$ws.Range("Z4").Value = 0.54
$ws.Range("AA4").Value = 768.3099999999999
$ws.Range("AB4").Value = 870.9299999999999
$ws.Range("AC4").Value = 2211
$ws.Range("AD4").Value = 13.91
$ws.Range("AE4").Value = 48552
$ws.Range("AF4").Value = 0.63
$ws.Range("AG4").Value = 630
$ws.Range("AH4").Value = 2.05
$ws.Range("AI4").Value = 28.51
$ws.Range("AJ4").Value = 78822179
$ws.Range("U4").ClearContents()

# Row 5
$ws.Range("D5").Value = 44855
$ws.Range("E5").Value = 3603
$ws.Range("F5").Value = 3603
$ws.Range("G5").Value = 3560
$ws.Range("H5").Value = 2710
$ws.Range("I5").Value = 2716
$ws.Range("J5").Value = -6
$ws.Range("K5").Value = 379483
$ws.Range("L5").Value = 335366
$ws.Range("M5").Value = 44116
$ws.Range("N5").Value = 44096
$ws.Range("O5").Value = 20
$ws.Range("P5").Value = 4585
$ws.Range("Q5").Value = -7891
$ws.Range("R5").Value = -11665
$ws.Range("S5").Value = 20891
$ws.Range("T5").Value = 164
$ws.Range("V5").Value = 17378
$ws.Range("W5").Value = 8.029999999999999
$ws.Range("X5").Value = 6.04
$ws.Range("Y5").Value = 6.6
$ws.Range("Z5").Value = 0.76
$ws.Range("AA5").Value = 760.1799999999999
$ws.Range("AB5").Value = 862.24
$ws.Range("AC5").Value = 3125
$ws.Range("AD5").Value = 11.71
$ws.Range("AE5").Value = 49380
$ws.Range("AF5").Value = 0.74
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 2.73
$ws.Range("AI5").Value = 32.88
$ws.Range("AJ5").Value = 89300000
$ws.Range("U5").ClearContents()

# Row 6
$ws.Range("D6").Value = 48902
$ws.Range("E6").Value = 4581
$ws.Range("F6").Value = 4581
$ws.Range("G6").Value = 4614
$ws.Range("H6").Value = 3341
$ws.Range("I6").Value = 3341
$ws.Range("K6").Value = 383306
$ws.Range("L6").Value = 336685
$ws.Range("M6").Value = 46621
$ws.Range("N6").Value = 46621
$ws.Range("P6").Value = 4585
$ws.Range("Q6").Value = -15202
$ws.Range("R6").Value = 27996
$ws.Range("S6").Value = -11967
$ws.Range("T6").Value = 114
$ws.Range("V6").Value = 26807
$ws.Range("W6").Value = 9.369999999999999
$ws.Range("X6").Value = 6.83
$ws.Range("Y6").Value = 7.36
$ws.Range("Z6").Value = 0.88
$ws.Range("AA6").Value = 722.17
$ws.Range("AB6").Value = 916.87
$ws.Range("AC6").Value = 3741
$ws.Range("AD6").Value = 8.42
$ws.Range("AE6").Value = 52207
$ws.Range("AF6").Value = 0.6
$ws.Range("AG6").Value = 1400
$ws.Range("AH6").Value = 4.44
$ws.Range("AI6").Value = 37.42
$ws.Range("AJ6").Value = 89300000
$ws.Range("U6").ClearContents()

# Row 7
$ws.Range("D7").Value = 65513
$ws.Range("E7").Value = 5035
$ws.Range("G7").Value = 5049
$ws.Range("H7").Value = 3753
$ws.Range("I7").Value = 3770
$ws.Range("K7").Value = 431242
$ws.Range("L7").Value = 382858
$ws.Range("M7").Value = 49228
$ws.Range("N7").Value = 48725
$ws.Range("P7").Value = 4583
$ws.Range("W7").Value = 7.69
$ws.Range("X7").Value = 5.73
$ws.Range("Y7").Value = 7.91
$ws.Range("Z7").Value = 0.92
$ws.Range("AA7").Value = 777.73
$ws.Range("AC7").Value = 4222
$ws.Range("AD7").Value = 9.02
$ws.Range("AE7").Value = 54563
$ws.Range("AF7").Value = 0.7
$ws.Range("AG7").Value = 1656
$ws.Range("AH7").Value = 4.35
$ws.Range("AI7").Value = 39.22
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()

# Row 8
$ws.Range("D8").Value = 57056
$ws.Range("E8").Value = 5039
$ws.Range("G8").Value = 5159
$ws.Range("H8").Value = 3801
$ws.Range("I8").Value = 3736
$ws.Range("K8").Value = 461010
$ws.Range("L8").Value = 409964
$ws.Range("M8").Value = 51734
$ws.Range("N8").Value = 50968
$ws.Range("P8").Value = 4582
$ws.Range("W8").Value = 8.83
$ws.Range("X8").Value = 6.66
$ws.Range("Y8").Value = 7.5
$ws.Range("Z8").Value = 0.85
$ws.Range("AA8").Value = 792.4400000000001
$ws.Range("AC8").Value = 4184
$ws.Range("AD8").Value = 8.43
$ws.Range("AE8").Value = 57075
$ws.Range("AF8").Value = 0.62
$ws.Range("AG8").Value = 1779
$ws.Range("AH8").Value = 5.05
$ws.Range("AI8").Value = 42.52
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()

# Row 9
$ws.Range("D9").Value = 63600
$ws.Range("E9").Value = 5412
$ws.Range("G9").Value = 5490
$ws.Range("H9").Value = 3979
$ws.Range("I9").Value = 4038
$ws.Range("K9").Value = 493044
$ws.Range("L9").Value = 438965
$ws.Range("M9").Value = 54079
$ws.Range("N9").Value = 54051
$ws.Range("P9").Value = 4581
$ws.Range("W9").Value = 8.51
$ws.Range("X9").Value = 6.26
$ws.Range("Y9").Value = 7.69
$ws.Range("Z9").Value = 0.83
$ws.Range("AA9").Value = 811.71
$ws.Range("AC9").Value = 4522
$ws.Range("AD9").Value = 7.8
$ws.Range("AE9").Value = 60528
$ws.Range("AF9").Value = 0.58
$ws.Range("AG9").Value = 1940
$ws.Range("AH9").Value = 5.5
$ws.Range("AI9").Value = 42.9
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
